$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "64.310.74"
$ws.Range("D3").Value = "3.502.13"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("D5").Value = "'590.90"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "'134.24"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").Value = "4.097.82"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "3.502.01"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "64.308.41"
$ws.Range("D17").Value = "'25.65"
$ws.Range("E17").Value = "  -6.75%  "
$ws.Range("D18").Value = "'9.85"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'5.76"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").Value = "'393.54"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "3.641.04"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").Value = "'74.59"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "'0.996"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "'7.38"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "'8.24"
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("E31").Value = "  -6.94%  "
$ws.Range("D32").Value = "3.523.28"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  +5.64%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'23.48"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("E36").Value = "  -5.50%  "
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").Value = "'167.34"
$ws.Range("E39").Value = "  +5.36%  "
$ws.Range("D40").Value = "'0.0780"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").Value = "'0.811"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D43").Value = "'25.16"
$ws.Range("E43").Value = "  -5.89%  "
$ws.Range("D44").Value = "'4.40"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("E46").Value = "  -3.72%  "
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "2.377.71"
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("D49").Value = "'0.894"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").Value = "'21.16"
$ws.Range("E51").Value = "  -1.56%  "

Write-Output "Applied cryptos update"
